$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Numeric-looking price strings are forced to remain text (matching the
# original inlineStr cell type) by applying a text number format before
# assignment, then resetting the cell style so no visible style change occurs.

$ws.Range("D2").Value = "69.163.90"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.473.39"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "2.471.26"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "69.040.60"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "2.467.27"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "339.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.71%  "
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").Value = "2.600.95"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").Value = "0.0₃0824"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "432.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -2.73%  "
$ws.Range("E35").Value = "  -3.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.06"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0716"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.563"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("E51").Value = "  +0.10%  "
